$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 800
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 800
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 800
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -1138
$ws.Range("H16").Value = 1700
$ws.Range("J16").Value = 1700
$ws.Range("L16").Value = 1700
$ws.Range("N16").Value = -2160
$ws.Range("H17").Value = 2356.739
$ws.Range("J17").Value = 2590.85
$ws.Range("L17").Value = 7772.549999999999
$ws.Range("N17").Value = -8108.549999999999
$ws.Range("H33").Value = 200
$ws.Range("J33").Value = 200
$ws.Range("L33").Value = 200
$ws.Range("N33").Value = -658
$ws.Range("H70").Value = 7074.75
$ws.Range("J70").Value = 8266.333000000001
$ws.Range("L70").Value = 24798.999
$ws.Range("N70").Value = -25338.999
$ws.Range("H73").Value = 7074.75
$ws.Range("J73").Value = 8266.333000000001
$ws.Range("L73").Value = 24798.999
$ws.Range("N73").Value = -26670.999
$ws.Range("H88").Value = 1485
$ws.Range("J88").Value = 975
$ws.Range("L88").Value = 975
$ws.Range("N88").Value = -1787
$ws.Range("H91").Value = 1485
$ws.Range("J91").Value = 975
$ws.Range("L91").Value = 975
$ws.Range("N91").Value = -3783
$ws.Range("H93").Value = 20000
$ws.Range("J93").Value = 20000
$ws.Range("L93").Value = 20000
$ws.Range("N93").Value = -24992
$ws.Range("H112").Value = 2232.2307
$ws.Range("J112").Value = 1946.5555
$ws.Range("L112").Value = 5839.666499999999
$ws.Range("N112").Value = -8055.666499999999
$ws.Range("H140").Value = 87666.664
$ws.Range("J140").Value = 87666.664
$ws.Range("L140").Value = 87666.664
$ws.Range("N140").Value = -98026.664

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3176
$ws.Range("J2").Value = 4287.625
$ws.Range("L2").Value = 4287.625
$ws.Range("N2").Value = -4513.625
$ws.Range("H32").Value = 3572144
$ws.Range("I32").Value = 630.52
$ws.Range("K32").Value = 630.52
$ws.Range("M32").Value = -343.52
$ws.Range("H45").Value = 2683.9443
$ws.Range("I45").Value = 2273.111
$ws.Range("J45").Value = 3094.7778
$ws.Range("K45").Value = 2273.111
$ws.Range("L45").Value = 3094.7778
$ws.Range("M45").Value = -1896.111
$ws.Range("N45").Value = -3848.7778
$ws.Range("H61").Value = 2785.5715
$ws.Range("I61").Value = 2300
$ws.Range("K61").Value = 2300
$ws.Range("M61").Value = -2088
$ws.Range("H82").Value = 30000
$ws.Range("J82").Value = 30000
$ws.Range("L82").Value = 30000
$ws.Range("N82").Value = -30722
$ws.Range("H85").Value = 30000
$ws.Range("J85").Value = 30000
$ws.Range("L85").Value = 30000
$ws.Range("N85").Value = -32496
$ws.Range("H116").Value = 3176
$ws.Range("J116").Value = 4287.625
$ws.Range("L116").Value = 4287.625
$ws.Range("N116").Value = -8875.625
$ws.Range("H136").Value = 2785.5715
$ws.Range("I136").Value = 2300
$ws.Range("K136").Value = 6900
$ws.Range("M136").Value = -4350

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3176
$ws.Range("J3").Value = 4287.625
$ws.Range("L3").Value = 4287.625
$ws.Range("N3").Value = -4515.625
$ws.Range("H94").Value = 243.33333
$ws.Range("I94").Value = 208.57143
$ws.Range("K94").Value = 208.57143
$ws.Range("M94").Value = 242.42857
$ws.Range("H107").Value = 33338234
$ws.Range("I107").Value = 125001500
$ws.Range("J107").Value = 6136.8184
$ws.Range("K107").Value = 125001500
$ws.Range("L107").Value = 6136.8184
$ws.Range("M107").Value = -124999580
$ws.Range("N107").Value = -9976.8184
$ws.Range("H134").Value = 8205.875
$ws.Range("J134").Value = 30199.5
$ws.Range("L134").Value = 90598.5
$ws.Range("N134").Value = -95668.5
$ws.Range("H140").Value = 60010
$ws.Range("J140").Value = 60010
$ws.Range("L140").Value = 60010
$ws.Range("N140").Value = -70370

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2816.6
$ws.Range("I22").Value = 1731.6666
$ws.Range("K22").Value = 1731.6666
$ws.Range("M22").Value = -1381.6666
$ws.Range("H31").Value = 5181.815
$ws.Range("I31").Value = 1863.875
$ws.Range("K31").Value = 1863.875
$ws.Range("M31").Value = -1568.875
$ws.Range("H34").Value = 5181.815
$ws.Range("I34").Value = 1863.875
$ws.Range("K34").Value = 1863.875
$ws.Range("M34").Value = -1661.875
$ws.Range("H88").Value = 6374
$ws.Range("J88").Value = 6374
$ws.Range("L88").Value = 6374
$ws.Range("N88").Value = -7186
$ws.Range("H91").Value = 6374
$ws.Range("J91").Value = 6374
$ws.Range("L91").Value = 6374
$ws.Range("N91").Value = -9182
$ws.Range("H105").Value = 3210.125
$ws.Range("I105").Value = 1417.5
$ws.Range("J105").Value = 5002.75
$ws.Range("K105").Value = 1417.5
$ws.Range("L105").Value = 5002.75
$ws.Range("M105").Value = 329.5
$ws.Range("N105").Value = -8496.75
$ws.Range("H132").Value = 3150.875
$ws.Range("J132").Value = 3997.5
$ws.Range("L132").Value = 11992.5
$ws.Range("N132").Value = -17052.5
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -55060
$ws.Range("H134").Value = 3082.2
$ws.Range("I134").Value = 3082.2
$ws.Range("K134").Value = 9246.599999999999
$ws.Range("M134").Value = -6711.599999999999
$ws.Range("H140").Value = 118750
$ws.Range("J140").Value = 118750
$ws.Range("L140").Value = 118750
$ws.Range("N140").Value = -129110

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 180
$ws.Range("I98").Value = 180
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 540
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 958
$ws.Range("N98").ClearContents()
$ws.Range("H137").Value = 4302.5
$ws.Range("J137").Value = 5758
$ws.Range("L137").Value = 17274
$ws.Range("N137").Value = -27474

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4903
$ws.Range("J70").Value = 6000
$ws.Range("L70").Value = 6000
$ws.Range("N70").Value = -6540
$ws.Range("H73").Value = 4903
$ws.Range("J73").Value = 6000
$ws.Range("L73").Value = 6000
$ws.Range("N73").Value = -7872
$ws.Range("H97").Value = 898.6667
$ws.Range("I97").Value = 200
$ws.Range("J97").Value = 1098.2858
$ws.Range("K97").Value = 200
$ws.Range("L97").Value = 1098.2858
$ws.Range("M97").Value = 296
$ws.Range("N97").Value = -2090.2858
$ws.Range("H102").Value = 4052.75
$ws.Range("I102").Value = 3703.6667
$ws.Range("K102").Value = 3703.6667
$ws.Range("M102").Value = -2081.6667
$ws.Range("H132").Value = 2003.375
$ws.Range("I132").Value = 2003.375
$ws.Range("K132").Value = 6010.125
$ws.Range("M132").Value = -3480.125

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2735.1428
$ws.Range("I40").Value = 3873
$ws.Range("J40").Value = 2280
$ws.Range("K40").Value = 3873
$ws.Range("L40").Value = 2280
$ws.Range("M40").Value = -3737
$ws.Range("N40").Value = -2552
$ws.Range("H46").Value = 5141.6665
$ws.Range("I46").Value = 1458.3334
$ws.Range("J46").Value = 6983.3335
$ws.Range("K46").Value = 1458.3334
$ws.Range("L46").Value = 6983.3335
$ws.Range("M46").Value = -1270.3334
$ws.Range("N46").Value = -7359.3335
$ws.Range("H82").Value = 3814.1428
$ws.Range("I82").Value = 625
$ws.Range("K82").Value = 625
$ws.Range("M82").Value = -264
$ws.Range("H85").Value = 3814.1428
$ws.Range("I85").Value = 625
$ws.Range("K85").Value = 625
$ws.Range("M85").Value = 623
$ws.Range("H106").Value = 7328.1665
$ws.Range("J106").Value = 7328.1665
$ws.Range("L106").Value = 7328.1665
$ws.Range("N106").Value = -9852.166499999999
$ws.Range("H122").Value = 3684.5715
$ws.Range("I122").Value = 3795.75
$ws.Range("K122").Value = 11387.25
$ws.Range("M122").Value = -8937.25
$ws.Range("H130").Value = 50000
$ws.Range("J130").Value = 50000
$ws.Range("L130").Value = 50000
$ws.Range("N130").Value = -60040
$ws.Range("H139").Value = 65000
$ws.Range("I139").Value = 65000
$ws.Range("K139").Value = 65000
$ws.Range("M139").Value = -59860

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 12300
$ws.Range("I52").Value = 30000
$ws.Range("J52").Value = 3450
$ws.Range("K52").Value = 30000
$ws.Range("L52").Value = 3450
$ws.Range("M52").Value = -29774
$ws.Range("N52").Value = -3902
$ws.Range("H81").Value = 776
$ws.Range("I81").Value = 776
$ws.Range("K81").Value = 1552
$ws.Range("M81").Value = -491
$ws.Range("H84").Value = 776
$ws.Range("I84").Value = 776
$ws.Range("K84").Value = 7760
$ws.Range("M84").Value = -2456
$ws.Range("H132").Value = 1059.2778
$ws.Range("I132").Value = 817.8
$ws.Range("K132").Value = 2453.4
$ws.Range("M132").Value = 76.60000000000036
